$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting (number formats/styles) from row 17 down into row 18
$ws.Range("A17:T17").Copy()
$ws.Range("A18:T18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row data: v16-4600 / Punishing going back
$ws.Range("A18").Value = "v16-4600"
$ws.Range("B18").Value = 3948
$ws.Range("C18").Value = "Punishing going back"
$ws.Range("D18").Value = 92
$ws.Range("E18").Value = 0.23
$ws.Range("F18").Value = 0.41
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 2461.96
$ws.Range("J18").Value = 1270
$ws.Range("K18").Value = 4130
$ws.Range("L18").Value = 213.89
$ws.Range("M18").Value = 119
$ws.Range("N18").Value = 244
$ws.Range("O18").Value = 0.78
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 4
$ws.Range("R18").Value = 110.75
$ws.Range("S18").Value = 49.7
$ws.Range("T18").Value = 215.3

# Update selection to match post-edit state
$ws.Range("C19").Select()
